# Apply "Mudanças na amostragem de dados cadastrados"
# Updates the sample schedule data and inserts a new class row ("2BADM")
# before the trailing header-like "Turma" row, pushing it from row 4 to row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (1BDS) : update D2, E2, F2 ---
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2

# --- Row 3 : rename 1ADS -> 3ADS and update values ---
$ws.Range("A3").Value = "3ADS"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1

# --- Move old row 4 ("Turma") down to row 5 ---
$ws.Range("A5").Value = "Turma"

# --- Row 4 becomes a new full data row: 2BADM ---
$ws.Range("A4").Value = "2BADM"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 3
